# Upload new version with timestamp
# Adds two new low-stock items (ALKAPRESS, STRINGAZOLE) to the report,
# renumbers the index column, recomputes the sell-price total, and
# refreshes the footer timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
}

# ---------------------------------------------------------------
# 1) Insert "ALKAPRESS PLUS 10/160MG 20 F.C. TABS." as the new row 7
#    (it sorts alphabetically first, ahead of AVIVAVASC).
# ---------------------------------------------------------------
$ws.Rows("7:7").Insert()
$ws.Range("A8:Q8").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)
$ws.Rows("7:7").RowHeight = 25.5
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()

Set-TextCell "C7" "ALKAPRESS PLUS 10/160MG 20 F.C. TABS."
Set-TextCell "H7" "0:1"
Set-TextCell "L7" "1"
Set-TextCell "N7" "102.00"
Set-TextCell "P7" "102.0000"
Set-TextCell "Q7" "1:0"

# ---------------------------------------------------------------
# 2) Insert "STRINGAZOLE 40MG 21 ENTERIC COATED TABLETS" as the new
#    row 16 (after the row-7 insert, SELGON sits at row 15 and TOBRIN
#    at row 16, so STRINGAZOLE belongs between them).
# ---------------------------------------------------------------
$ws.Rows("16:16").Insert()
$ws.Range("A15:Q15").Copy()
$ws.Range("A16:Q16").PasteSpecial(-4122)
$ws.Rows("16:16").RowHeight = 25.5
$ws.Range("A16:B16").Merge()
$ws.Range("C16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()
$ws.Range("N16:O16").Merge()

Set-TextCell "C16" "STRINGAZOLE 40MG 21 ENTERIC COATED TABLETS"
Set-TextCell "H16" "1:1"
Set-TextCell "L16" "1"
Set-TextCell "N16" "126.00"
Set-TextCell "P16" "41.5800"
Set-TextCell "Q16" "0:1"

# ---------------------------------------------------------------
# 3) Renumber the "م" index column (A7:A17) sequentially 1..11.
# ---------------------------------------------------------------
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 9
$ws.Range("A16").Value = 10
$ws.Range("A17").Value = 11

# ---------------------------------------------------------------
# 4) Update the sell-price total (sum of column P across all items).
# ---------------------------------------------------------------
Set-TextCell "P18" "682.08000000000004"
$ws.Range("P18").NumberFormat = "#.00"
$ws.Range("P18").Value = 682.08000000000004

# ---------------------------------------------------------------
# 5) Refresh the generated-at timestamp in the footer.
# ---------------------------------------------------------------
$ws.Range("A19").Value = "Thursday, 10 July, 2025 10:46 AM"

Write-Output "edit complete"
